# Applies the commit "[EXTRA SCRAPE] full data scraped for extra batting and
# bowling fields" to the "ODI Batting Extra" worksheet (4th sheet / sheetId 4):
# the used range grows from A1:F21 to A1:F76, with row 2 onward re-populated from
# the newly scraped MATCH_CODE/BATTING_POSITION/NUM_4/NUM_6/PERCENT_RUNS_OF_TOTAL/
# MAN_OF_MATCH data (75 data rows total).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Columns A, C, D, E, F hold text (match codes / flags / formatted percents),
# even when a value looks numeric (e.g. "0", "1", "1.55%") - force text format
# first so Excel does not auto-convert them to numbers/percentages.
$ws.Range("A2:A76").NumberFormat = "@"
$ws.Range("C2:F76").NumberFormat = "@"
# Column B (BATTING_POSITION) holds real numbers.
$ws.Range("B2:B76").NumberFormat = "General"

# Row 2
$ws.Cells.Item(2, 1).Value = "3111"
$ws.Cells.Item(2, 6).Value = "NO"

# Row 3
$ws.Cells.Item(3, 1).Value = "3113"
$ws.Cells.Item(3, 2).Value = 11
$ws.Cells.Item(3, 3).Value = "0"
$ws.Cells.Item(3, 4).Value = "0"
$ws.Cells.Item(3, 5).Value = "1.55%"
$ws.Cells.Item(3, 6).Value = "NO"

# Row 4
$ws.Cells.Item(4, 1).Value = "3116"
$ws.Cells.Item(4, 2).Value = 11
$ws.Cells.Item(4, 6).Value = "NO"

# Row 5
$ws.Cells.Item(5, 1).Value = "3335"
$ws.Cells.Item(5, 2).Value = 10
$ws.Cells.Item(5, 3).Value = ""
$ws.Cells.Item(5, 4).Value = ""
$ws.Cells.Item(5, 5).Value = ""
$ws.Cells.Item(5, 6).Value = "NO"

# Row 6
$ws.Cells.Item(6, 1).Value = "3337"
$ws.Cells.Item(6, 2).Value = ""
$ws.Cells.Item(6, 6).Value = "NO"

# Row 7
$ws.Cells.Item(7, 1).Value = "3341"
$ws.Cells.Item(7, 2).Value = ""
$ws.Cells.Item(7, 6).Value = "NO"

# Row 8
$ws.Cells.Item(8, 1).Value = "3353"
$ws.Cells.Item(8, 2).Value = ""
$ws.Cells.Item(8, 6).Value = "NO"

# Row 9
$ws.Cells.Item(9, 1).Value = "3355"
$ws.Cells.Item(9, 2).Value = 10
$ws.Cells.Item(9, 6).Value = "NO"

# Row 10
$ws.Cells.Item(10, 1).Value = "3357"
$ws.Cells.Item(10, 2).Value = 11
$ws.Cells.Item(10, 3).Value = "0"
$ws.Cells.Item(10, 4).Value = "0"
$ws.Cells.Item(10, 6).Value = "NO"

# Row 11
$ws.Cells.Item(11, 1).Value = "3373"
$ws.Cells.Item(11, 2).Value = 11
$ws.Cells.Item(11, 6).Value = "NO"

# Row 12
$ws.Cells.Item(12, 1).Value = "3375"
$ws.Cells.Item(12, 6).Value = "NO"

# Row 13
$ws.Cells.Item(13, 1).Value = "3380"
$ws.Cells.Item(13, 2).Value = ""
$ws.Cells.Item(13, 3).Value = ""
$ws.Cells.Item(13, 4).Value = ""
$ws.Cells.Item(13, 6).Value = "NO"

# Row 14
$ws.Cells.Item(14, 1).Value = "3382"
$ws.Cells.Item(14, 2).Value = 11
$ws.Cells.Item(14, 6).Value = "NO"

# Row 15
$ws.Cells.Item(15, 1).Value = "3386"
$ws.Cells.Item(15, 2).Value = 11
$ws.Cells.Item(15, 6).Value = "NO"

# Row 16
$ws.Cells.Item(16, 1).Value = "3387"
$ws.Cells.Item(16, 2).Value = 11
$ws.Cells.Item(16, 6).Value = "NO"

# Row 17
$ws.Cells.Item(17, 1).Value = "3431"
$ws.Cells.Item(17, 2).Value = 10
$ws.Cells.Item(17, 3).Value = "0"
$ws.Cells.Item(17, 4).Value = "0"
$ws.Cells.Item(17, 5).Value = "0.68%"
$ws.Cells.Item(17, 6).Value = "NO"

# Row 18
$ws.Cells.Item(18, 1).Value = "3432"
$ws.Cells.Item(18, 6).Value = "NO"

# Row 19
$ws.Cells.Item(19, 1).Value = "3505"
$ws.Cells.Item(19, 2).Value = 10
$ws.Cells.Item(19, 6).Value = "NO"

# Row 20
$ws.Cells.Item(20, 1).Value = "3510"
$ws.Cells.Item(20, 2).Value = 11
$ws.Cells.Item(20, 3).Value = "2"
$ws.Cells.Item(20, 4).Value = "1"
$ws.Cells.Item(20, 5).Value = "16.19%"
$ws.Cells.Item(20, 6).Value = "NO"

# Row 21
$ws.Cells.Item(21, 1).Value = "3514"
$ws.Cells.Item(21, 6).Value = "NO"

# Row 22
$ws.Cells.Item(22, 1).Value = "3518"
$ws.Cells.Item(22, 2).Value = 11
$ws.Cells.Item(22, 6).Value = "NO"

# Row 23
$ws.Cells.Item(23, 1).Value = "3519"
$ws.Cells.Item(23, 2).Value = 11
$ws.Cells.Item(23, 6).Value = "NO"

# Row 24
$ws.Cells.Item(24, 1).Value = "3521"
$ws.Cells.Item(24, 2).Value = 9
$ws.Cells.Item(24, 6).Value = "NO"

# Row 25
$ws.Cells.Item(25, 1).Value = "3524"
$ws.Cells.Item(25, 6).Value = "NO"

# Row 26
$ws.Cells.Item(26, 1).Value = "3525"
$ws.Cells.Item(26, 6).Value = "NO"

# Row 27
$ws.Cells.Item(27, 1).Value = "3529"
$ws.Cells.Item(27, 2).Value = 10
$ws.Cells.Item(27, 6).Value = "NO"

# Row 28
$ws.Cells.Item(28, 1).Value = "3587"
$ws.Cells.Item(28, 6).Value = "NO"

# Row 29
$ws.Cells.Item(29, 1).Value = "3588"
$ws.Cells.Item(29, 2).Value = 11
$ws.Cells.Item(29, 3).Value = "0"
$ws.Cells.Item(29, 4).Value = "0"
$ws.Cells.Item(29, 6).Value = "NO"

# Row 30
$ws.Cells.Item(30, 1).Value = "3643"
$ws.Cells.Item(30, 2).Value = 10
$ws.Cells.Item(30, 6).Value = "NO"

# Row 31
$ws.Cells.Item(31, 1).Value = "3644"
$ws.Cells.Item(31, 2).Value = 10
$ws.Cells.Item(31, 6).Value = "NO"

# Row 32
$ws.Cells.Item(32, 1).Value = "3645"
$ws.Cells.Item(32, 2).Value = 10
$ws.Cells.Item(32, 6).Value = "NO"

# Row 33
$ws.Cells.Item(33, 1).Value = "3672"
$ws.Cells.Item(33, 6).Value = "NO"

# Row 34
$ws.Cells.Item(34, 1).Value = "3680"
$ws.Cells.Item(34, 2).Value = 11
$ws.Cells.Item(34, 6).Value = "NO"

# Row 35
$ws.Cells.Item(35, 1).Value = "3683"
$ws.Cells.Item(35, 2).Value = 11
$ws.Cells.Item(35, 6).Value = "NO"

# Row 36
$ws.Cells.Item(36, 1).Value = "3688"
$ws.Cells.Item(36, 2).Value = 10
$ws.Cells.Item(36, 3).Value = "1"
$ws.Cells.Item(36, 4).Value = "0"
$ws.Cells.Item(36, 5).Value = "1.62%"
$ws.Cells.Item(36, 6).Value = "NO"

# Row 37
$ws.Cells.Item(37, 1).Value = "3689"
$ws.Cells.Item(37, 2).Value = 9
$ws.Cells.Item(37, 6).Value = "NO"

# Row 38
$ws.Cells.Item(38, 1).Value = "3692"
$ws.Cells.Item(38, 6).Value = "NO"

# Row 39
$ws.Cells.Item(39, 1).Value = "3693"
$ws.Cells.Item(39, 2).Value = 10
$ws.Cells.Item(39, 6).Value = "NO"

# Row 40
$ws.Cells.Item(40, 1).Value = "3731"
$ws.Cells.Item(40, 6).Value = "NO"

# Row 41
$ws.Cells.Item(41, 1).Value = "3735"
$ws.Cells.Item(41, 6).Value = "NO"

# Row 42
$ws.Cells.Item(42, 1).Value = "3751"
$ws.Cells.Item(42, 6).Value = "NO"

# Row 43
$ws.Cells.Item(43, 1).Value = "3760"
$ws.Cells.Item(43, 2).Value = 10
$ws.Cells.Item(43, 3).Value = "1"
$ws.Cells.Item(43, 4).Value = "0"
$ws.Cells.Item(43, 5).Value = "7.63%"
$ws.Cells.Item(43, 6).Value = "NO"

# Row 44
$ws.Cells.Item(44, 1).Value = "3768"
$ws.Cells.Item(44, 2).Value = 10
$ws.Cells.Item(44, 6).Value = "NO"

# Row 45
$ws.Cells.Item(45, 1).Value = "3775"
$ws.Cells.Item(45, 2).Value = 10
$ws.Cells.Item(45, 6).Value = "NO"

# Row 46
$ws.Cells.Item(46, 1).Value = "3781"
$ws.Cells.Item(46, 2).Value = 10
$ws.Cells.Item(46, 6).Value = "NO"

# Row 47
$ws.Cells.Item(47, 1).Value = "3786"
$ws.Cells.Item(47, 2).Value = 10
$ws.Cells.Item(47, 6).Value = "NO"

# Row 48
$ws.Cells.Item(48, 1).Value = "3791"
$ws.Cells.Item(48, 6).Value = "NO"

# Row 49
$ws.Cells.Item(49, 1).Value = "3795"
$ws.Cells.Item(49, 2).Value = 11
$ws.Cells.Item(49, 6).Value = "NO"

# Row 50
$ws.Cells.Item(50, 1).Value = "3808"
$ws.Cells.Item(50, 6).Value = "NO"

# Row 51
$ws.Cells.Item(51, 1).Value = "3811"
$ws.Cells.Item(51, 2).Value = 10
$ws.Cells.Item(51, 3).Value = "0"
$ws.Cells.Item(51, 4).Value = "0"
$ws.Cells.Item(51, 6).Value = "NO"

# Row 52
$ws.Cells.Item(52, 1).Value = "3841"
$ws.Cells.Item(52, 2).Value = 10
$ws.Cells.Item(52, 6).Value = "NO"

# Row 53
$ws.Cells.Item(53, 1).Value = "3844"
$ws.Cells.Item(53, 6).Value = "NO"

# Row 54
$ws.Cells.Item(54, 1).Value = "3874"
$ws.Cells.Item(54, 2).Value = 10
$ws.Cells.Item(54, 6).Value = "NO"

# Row 55
$ws.Cells.Item(55, 1).Value = "3875"
$ws.Cells.Item(55, 2).Value = 8
$ws.Cells.Item(55, 3).Value = "0"
$ws.Cells.Item(55, 4).Value = "0"
$ws.Cells.Item(55, 6).Value = "NO"

# Row 56
$ws.Cells.Item(56, 1).Value = "3876"
$ws.Cells.Item(56, 6).Value = "NO"

# Row 57
$ws.Cells.Item(57, 1).Value = "3877"
$ws.Cells.Item(57, 2).Value = 10
$ws.Cells.Item(57, 6).Value = "NO"

# Row 58
$ws.Cells.Item(58, 1).Value = "3878"

# Row 59
$ws.Cells.Item(59, 1).Value = "3951"

# Row 60
$ws.Cells.Item(60, 1).Value = "3952"

# Row 61
$ws.Cells.Item(61, 1).Value = "3953"

# Row 62
$ws.Cells.Item(62, 1).Value = "3954"

# Row 63
$ws.Cells.Item(63, 1).Value = "3955"

# Row 64
$ws.Cells.Item(64, 1).Value = "3974"

# Row 65
$ws.Cells.Item(65, 1).Value = "4034"

# Row 66
$ws.Cells.Item(66, 1).Value = "4038"

# Row 67
$ws.Cells.Item(67, 1).Value = "4051"

# Row 68
$ws.Cells.Item(68, 1).Value = "4052"

# Row 69
$ws.Cells.Item(69, 1).Value = "4053"

# Row 70
$ws.Cells.Item(70, 1).Value = "4056"

# Row 71
$ws.Cells.Item(71, 1).Value = "4057"

# Row 72
$ws.Cells.Item(72, 1).Value = "4074"

# Row 73
$ws.Cells.Item(73, 1).Value = "4171"

# Row 74
$ws.Cells.Item(74, 1).Value = "4173"

# Row 75
$ws.Cells.Item(75, 1).Value = "4213"

# Row 76
$ws.Cells.Item(76, 1).Value = "4216"

Write-Output "ODI Batting Extra updated: rows 2-76 (A:F) written"